# ValueSet-FrMethodOfAdministration.xlsx update
#  - refresh URL / Date / Jurisdiction / Description metadata
#  - swap the SNOMED CT expression-constraint include for an EDQM "All codes" include

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Include #0")

# ---- Metadata sheet -------------------------------------------------------
$ws1.Range("B2").Value  = "https://hl7.fr/ig/fhir/medication/ValueSet/FrMethodOfAdministration"
$ws1.Range("B8").Value  = "2026-01-15T08:54:26+00:00"
$ws1.Range("B11").Value = "FRANCE"
$ws1.Range("B12").Value = "Le jeu de valeurs à utiliser pour coder l'élément *dosageInstruction.method* de la ressource *FRMedicationRequest*."

# ---- Include #0 sheet ------------------------------------------------------
# Drop column C (the old "Value" header / constraint-expression column) entirely.
$ws2.Columns.Item(3).Delete()

# Drop the old Operation/constraint row content in A1:B2 (fully removed, not
# just blanked) - will be replaced with the new Codes/All-codes rows below.
$ws2.Range("B1:B2").Clear()

$ws2.Range("A1").Value = "Codes"
$ws2.Range("A2").Value = "All codes"

# Row 3 keeps its blank/blank cells (A3 & B3 still point at the shared blank
# string) - nothing to change there content-wise.

# Row 4: System URI now points at EDQM's Standard Terms instead of SNOMED CT.
$ws2.Range("A4").Value = "System URI"
$ws2.Range("B4").Value = "http://standardterms.edqm.eu"

# Re-apply the normal row style to the new B3/B4 cells (PasteSpecial formats
# only, so it reuses the existing style index instead of minting a new one).
$ws2.Range("A3").Copy()
$ws2.Range("B3").PasteSpecial(-4122)
$ws2.Range("A4").Copy()
$ws2.Range("B4").PasteSpecial(-4122)

$excel.CutCopyMode = 0
